# Apply updated "想去人数" (F column) / "最低票价" (G column) counts scraped
# at commit 456a3b4 to the 北京-漫展信息 workbook.
#
# Sheet "展览" (Exhibition) and sheet "全部类型" (All types) share the same
# rows/values (全部类型 aggregates 展览 + 演出 + 本地生活), so both need the
# same F-column updates. Sheet "演出" (Performance) has its own small set of
# updates, including a G-column value. Sheet "本地生活" (Local life) is
# unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: 展览 (Exhibition)
# ---------------------------------------------------------------------
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 11
$wsExhibition.Range("F3").Value = 103
$wsExhibition.Range("F4").Value = 933
$wsExhibition.Range("F5").Value = 66
$wsExhibition.Range("F6").Value = 7239
$wsExhibition.Range("F9").Value = 6594
$wsExhibition.Range("F11").Value = 280
$wsExhibition.Range("F12").Value = 4569
$wsExhibition.Range("F16").Value = 4666
$wsExhibition.Range("F17").Value = 20
$wsExhibition.Range("F20").Value = 355
$wsExhibition.Range("F28").Value = 8286
$wsExhibition.Range("F30").Value = 1446
$wsExhibition.Range("F31").Value = 62
$wsExhibition.Range("F37").Value = 1706
$wsExhibition.Range("F38").Value = 217
$wsExhibition.Range("F39").Value = 984
$wsExhibition.Range("F41").Value = 4322
$wsExhibition.Range("F46").Value = 857
$wsExhibition.Range("F47").Value = 1141

# ---------------------------------------------------------------------
# Sheet: 演出 (Performance)
# ---------------------------------------------------------------------
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 25
$wsPerformance.Range("F5").Value = 23
$wsPerformance.Range("F8").Value = 16
$wsPerformance.Range("G21").Value = 153

# ---------------------------------------------------------------------
# Sheet: 全部类型 (All types)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 25
$wsAll.Range("F5").Value = 11
$wsAll.Range("F6").Value = 103
$wsAll.Range("F7").Value = 23
$wsAll.Range("F8").Value = 933
$wsAll.Range("F9").Value = 66
$wsAll.Range("F10").Value = 7239
$wsAll.Range("F13").Value = 6594
$wsAll.Range("F15").Value = 280
$wsAll.Range("F16").Value = 4569
$wsAll.Range("F20").Value = 4666
$wsAll.Range("F21").Value = 20
$wsAll.Range("F23").Value = 355
$wsAll.Range("F29").Value = 8286
$wsAll.Range("F31").Value = 1446
$wsAll.Range("F32").Value = 62
$wsAll.Range("F37").Value = 1706
$wsAll.Range("F38").Value = 217
$wsAll.Range("F39").Value = 984
$wsAll.Range("F41").Value = 4322
$wsAll.Range("F46").Value = 857
$wsAll.Range("F47").Value = 1141
